$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# --- 1. Set the EARNED value for the 6/1/2023 period row (row 75) ---------
$ws.Range("C75").Value = 1.25

# --- 2. Insert a new table row above row 76 (new SICK LEAVE entry) --------
# A plain worksheet row insert shifts every cell in the row (incl. the
# stray out-of-table cell in column M), so do that and then drop the
# leftover M-column cell that results (it isn't part of Table1, and the
# author's edit dropped it rather than shifting it down with everything
# else).
$ws.Range("A76").EntireRow.Insert()
$lo.Resize($ws.Range("A8:K143"))
$ws.Range("M78").ClearContents()

# A bare row insert leaves the new row with generic (borderless) styles,
# so pull the real table-row formatting back from row 77 (the row that
# used to be row 76, carrying the correct bordered style set).
$ws.Range("A77:K77").Copy()
$ws.Range("A76:K76").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Fill in the freshly inserted row (matches the SL(1-0-0) pattern
#        already used for similar entries, e.g. row 69) -------------------
$ws.Range("B76").Value = "SL(1-0-0)"
$ws.Range("H76").Value = 1
$ws.Range("K76").Value = 45107

# Pull the date-number-format style straight off the analogous K69 cell
# (setting a NumberFormat string from scratch mints a brand-new style
# instead of reusing the workbook's existing date-format style).
$ws.Range("K69").Copy()
$ws.Range("K76").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Formats-only paste doesn't carry the calculated-column formula, so
# restore it explicitly on the new row.
$ws.Range("G76").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- 4. Re-assert the calculated-column formula text for the new last
#        table row (row 143) so it serialises the same way as its peers --
$ws.Range("G143").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- 5. Restore the active selection to the cell the author was last on --
$ws.Range("K76").Select()

$excel.CalculateFull()
